$d = $word.ActiveDocument

# 1. Fix text: remove "can" from "How can " -> "How "
$d.Content.Find.Execute("Why is it important to annotate our data? How can ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Why is it important to annotate our data? How ", 2)
